$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6:A35").NumberFormat = "General"
$ws.Range("A6:A35").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
